# Update the cryptos price/volume table (columns D = Price, E = Volume(1h))
# with the latest scraped values. Numeric-looking Price values are written
# with a leading apostrophe so Excel stores them as text (preserving exact
# formatting such as trailing zeros, e.g. "101.30", "5.542", "10.02")
# instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '31.385.14'
$ws.Range("E2").Value = '  +3.31%  '
$ws.Range("D3").Value = '1.996.67'
$ws.Range("E3").Value = '  +6.69%  '
$ws.Range("D4").Value = '''0.9989'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '''0.8093'
$ws.Range("E5").Value = '  +71.75%  '
$ws.Range("D6").Value = '''254.73'
$ws.Range("E6").Value = '  +4.41%  '
$ws.Range("D7").Value = '''0.9987'
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").Value = '''0.3524'
$ws.Range("E8").Value = '  +22.11%  '
$ws.Range("D9").Value = '''25.62'
$ws.Range("E9").Value = '  +16.54%  '
$ws.Range("D10").Value = '''0.07026'
$ws.Range("E10").Value = '  +8.59%  '
$ws.Range("D11").Value = '''0.8419'
$ws.Range("E11").Value = '  +16.19%  '
$ws.Range("D12").Value = '''0.08118'
$ws.Range("E12").Value = '  +4.18%  '
$ws.Range("D13").Value = '''101.30'
$ws.Range("E13").Value = '  +5.28%  '
$ws.Range("D14").Value = '1.989.68'
$ws.Range("E14").Value = '  +6.36%  '
$ws.Range("D15").Value = '''5.542'
$ws.Range("E15").Value = '  +7.79%  '
$ws.Range("D16").Value = '''272.75'
$ws.Range("E16").Value = '  -3.38%  '
$ws.Range("D17").Value = '31.357.75'
$ws.Range("E17").Value = '  +3.26%  '
$ws.Range("E18").Value = '  +7.49%  '
$ws.Range("D19").Value = '''0.000007947'
$ws.Range("E19").Value = '  +5.81%  '
$ws.Range("D20").Value = '''5.830'
$ws.Range("E20").Value = '  +10.78%  '
$ws.Range("D21").Value = '2.252.53'
$ws.Range("E21").Value = '  +6.75%  '
$ws.Range("D22").Value = '''0.9989'
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").Value = '''0.9990'
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").Value = '''6.977'
$ws.Range("E24").Value = '  +11.71%  '
$ws.Range("E25").Value = '  +8.42%  '
$ws.Range("D26").Value = '''0.1514'
$ws.Range("E26").Value = '  +57.36%  '
$ws.Range("D27").Value = '''164.51'
$ws.Range("E27").Value = '  +0.52%  '
$ws.Range("D28").Value = '''20.09'
$ws.Range("E28").Value = '  +7.15%  '
$ws.Range("D29").Value = '''2.233'
$ws.Range("E29").Value = '  +18.61%  '
$ws.Range("D30").Value = '''1.573'
$ws.Range("E30").Value = '  +5.83%  '
$ws.Range("D31").Value = '''1.358'
$ws.Range("E31").Value = '  +2.83%  '
$ws.Range("D32").Value = '''4.600'
$ws.Range("E32").Value = '  +8.63%  '
$ws.Range("D33").Value = '''4.353'
$ws.Range("E33").Value = '  +5.71%  '
$ws.Range("D34").Value = '''0.05207'
$ws.Range("E34").Value = '  +8.00%  '
$ws.Range("E35").Value = '  +8.39%  '
$ws.Range("D36").Value = '''0.7601'
$ws.Range("E36").Value = '  +10.10%  '
$ws.Range("D37").Value = '''2.773'
$ws.Range("E37").Value = '  +2.11%  '
$ws.Range("D38").Value = '''0.02011'
$ws.Range("E38").Value = '  +6.29%  '
$ws.Range("D39").Value = '''2.911'
$ws.Range("E39").Value = '  +3.21%  '
$ws.Range("D40").Value = '''6.652'
$ws.Range("E40").Value = '  +6.51%  '
$ws.Range("D41").Value = '''0.4758'
$ws.Range("E41").Value = '  +12.42%  '
$ws.Range("D42").Value = '''78.51'
$ws.Range("E42").Value = '  +4.05%  '
$ws.Range("D43").Value = '''2.124'
$ws.Range("E43").Value = '  +9.89%  '
$ws.Range("D44").Value = '''0.8616'
$ws.Range("E44").Value = '  +4.05%  '
$ws.Range("D45").Value = '''104.38'
$ws.Range("E45").Value = '  +3.41%  '
$ws.Range("D46").Value = '''0.9993'
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").Value = '''10.02'
$ws.Range("E47").Value = '  +3.59%  '
$ws.Range("D48").Value = '''7.533'
$ws.Range("E48").Value = '  +8.01%  '
$ws.Range("D49").Value = '''0.4384'
$ws.Range("E49").Value = '  +11.47%  '
$ws.Range("D50").Value = '''36.95'
$ws.Range("E50").Value = '  +4.60%  '
$ws.Range("D51").Value = '''0.1206'
$ws.Range("E51").Value = '  +13.66%  '
